# Updates crypto price/volume values per the latest data refresh.
# Uses Set-CellValue so that numeric-looking strings (e.g. "1.00", "584.16")
# are preserved as literal text, matching the source data's text formatting,
# instead of being auto-converted to numbers by Excel.

function Set-CellValue($Ws, $Cell, $Value) {
    $range = $Ws.Range($Cell)
    if ($Value -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number - force text storage so it keeps its
        # original textual representation (e.g. "1.00" instead of 1).
        $range.NumberFormat = "@"
    }
    $range.Value = $Value
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellValue $ws "D2" "67.089.50"
Set-CellValue $ws "E2" "  -0.40%  "
Set-CellValue $ws "D3" "2.477.67"
Set-CellValue $ws "E3" "  -1.34%  "
Set-CellValue $ws "D4" "1.00"
Set-CellValue $ws "E4" "  -0.06%  "
Set-CellValue $ws "D5" "584.16"
Set-CellValue $ws "E5" "  -0.77%  "
Set-CellValue $ws "D6" "169.55"
Set-CellValue $ws "E6" "  +0.04%  "
Set-CellValue $ws "E7" "  -0.02%  "
Set-CellValue $ws "E8" "  -1.54%  "
Set-CellValue $ws "D9" "2.478.23"
Set-CellValue $ws "E9" "  -1.17%  "
Set-CellValue $ws "D10" "0.138"
Set-CellValue $ws "E10" "  +1.46%  "
Set-CellValue $ws "E11" "  -1.18%  "
Set-CellValue $ws "D12" "4.92"
Set-CellValue $ws "E12" "  -1.00%  "
Set-CellValue $ws "E13" "  -2.92%  "
Set-CellValue $ws "D14" "2.928.09"
Set-CellValue $ws "E14" "  -2.44%  "
Set-CellValue $ws "D15" "25.26"
Set-CellValue $ws "E15" "  -2.84%  "
Set-CellValue $ws "D16" "66.859.29"
Set-CellValue $ws "E16" "  -0.72%  "
Set-CellValue $ws "D17" "0.0000170"
Set-CellValue $ws "E17" "  -2.09%  "
Set-CellValue $ws "D18" "2.541.68"
Set-CellValue $ws "E18" "  +0.15%  "
Set-CellValue $ws "D19" "11.00"
Set-CellValue $ws "E19" "  -5.95%  "
Set-CellValue $ws "E20" "  -7.39%  "
Set-CellValue $ws "D21" "350.41"
Set-CellValue $ws "E21" "  -3.79%  "
Set-CellValue $ws "E22" "  -1.71%  "
Set-CellValue $ws "E23" "  +0.29%  "
Set-CellValue $ws "D24" "68.50"
Set-CellValue $ws "E24" "  -4.54%  "
Set-CellValue $ws "E25" "  -5.90%  "
Set-CellValue $ws "E26" "  -2.24%  "
Set-CellValue $ws "D27" "9.13"
Set-CellValue $ws "E27" "  -5.64%  "
Set-CellValue $ws "D28" "0.998"
Set-CellValue $ws "E28" "  -27.23%  "
Set-CellValue $ws "D29" "2.587.91"
Set-CellValue $ws "E29" "  -2.80%  "
Set-CellValue $ws "E30" "  -4.14%  "
Set-CellValue $ws "D31" "507.97"
Set-CellValue $ws "E31" "  -3.80%  "
Set-CellValue $ws "D32" "7.65"
Set-CellValue $ws "E32" "  -6.72%  "
Set-CellValue $ws "D33" "1.76"
Set-CellValue $ws "E33" "  -4.36%  "
Set-CellValue $ws "E34" "  -3.62%  "
Set-CellValue $ws "D35" "1.00"
Set-CellValue $ws "E35" "  -0.03%  "
Set-CellValue $ws "E36" "  +0.64%  "
Set-CellValue $ws "D37" "0.116"
Set-CellValue $ws "E37" "  -9.13%  "
Set-CellValue $ws "D39" "18.22"
Set-CellValue $ws "E39" "  -4.79%  "
Set-CellValue $ws "E40" "  -6.51%  "
Set-CellValue $ws "E41" "  -0.15%  "
Set-CellValue $ws "E42" "  -3.86%  "
Set-CellValue $ws "E43" "  -3.87%  "
Set-CellValue $ws "E44" "  -4.14%  "
Set-CellValue $ws "E45" "  -2.66%  "
Set-CellValue $ws "D46" "38.84"
Set-CellValue $ws "E46" "  -1.24%  "
Set-CellValue $ws "D47" "142.09"
Set-CellValue $ws "E47" "  -3.20%  "
Set-CellValue $ws "B48" "Filecoin"
Set-CellValue $ws "C48" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellValue $ws "D48" "3.44"
Set-CellValue $ws "E48" "  -6.58%  "
Set-CellValue $ws "B49" "ARBITRUM"
Set-CellValue $ws "C49" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellValue $ws "D49" "0.512"
Set-CellValue $ws "E49" "  -5.24%  "
Set-CellValue $ws "E50" "  -6.89%  "
Set-CellValue $ws "D51" "0.0729"
Set-CellValue $ws "E51" "  -1.42%  "
